$d = $word.ActiveDocument

# Locate the paragraph that currently reads
# "GERENCIA DE LICENCIAS Y DESARROLLO ECONÓMICO" (built up out of five
# separate runs) and replace its contents with a single run that reads
# "SUBGERENCIA DE LICENCIAS COMERCIALES Y AUTORIZACIONES", while also
# dropping the paragraph's explicit "space after" override.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*GERENCIA DE LICENCIAS*DESARROLLO*") {
        $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Sinespaciado"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>SUBGERENCIA DE LICENCIAS COMERCIALES Y AUTORIZACIONES</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $p.Range.InsertXML($xml)
        break
    }
}
